$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("Lido2021")

$text = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."
$ws.Range("A41").Value = $text

$ws.Range("D2").Value = 0.03168685327642729
$ws.Range("E2").Value = -0.002842866988283843
$ws.Range("D3").Value = 0.02828851859495853
$ws.Range("E3").Value = 0.01058704327867499
$ws.Range("D4").Value = 0.02885552515033814
$ws.Range("E4").Value = 0.01181152907947114
$ws.Range("D5").Value = 0.06446813839444528
$ws.Range("E5").Value = -0.004464110187305081
$ws.Range("D6").Value = 0.01592570888249501
$ws.Range("E6").Value = -0.004627929185337498
$ws.Range("D7").Value = 0.01591050031601646
$ws.Range("E7").Value = 0.01360294117647043
$ws.Range("D8").Value = 0.0293993288927569
$ws.Range("E8").Value = 0.01157978511738955
$ws.Range("D9").Value = 0.03444545325768956
$ws.Range("E9").Value = 0.01935922110268318
$ws.Range("D10").Value = 0.02922092070906626
$ws.Range("E10").Value = 0.008587728956060303
$ws.Range("D11").Value = 0.0309034171211607
$ws.Range("E11").Value = -0.004902393781467951
$ws.Range("D12").Value = 0.01112760114013553
$ws.Range("E12").Value = 0.04608375678990706
$ws.Range("D13").Value = 0.01442669017625071
$ws.Range("E13").Value = 0.01929990539262061
$ws.Range("D14").Value = 0.01414689154936979
$ws.Range("E14").Value = 0.03685479980704276
$ws.Range("D15").Value = 0.009182074520610874
$ws.Range("E15").Value = 0.01469464027860345
$ws.Range("D16").Value = 0.007879207325615504
$ws.Range("E16").Value = 0.02227171492204905
$ws.Range("D17").Value = 0.02976101980067614
$ws.Range("E17").Value = -0.007108461362072993
$ws.Range("D18").Value = 0.02610004487502019
$ws.Range("E18").Value = -0.0054236173884461
$ws.Range("D19").Value = 0.03222461258857893
$ws.Range("E19").Value = 0.00571791613723005
$ws.Range("D20").Value = 0.03119900925938473
$ws.Range("E20").Value = -0.002937316417723879
$ws.Range("D21").Value = 0.04643194844061985
$ws.Range("E21").Value = 0.007281583975476069
$ws.Range("D22").Value = 0.03567871201380236
$ws.Range("E22").Value = 0.01344372489548307
$ws.Range("D23").Value = 0.03291094789632866
$ws.Range("E23").Value = 0.004946975531725961
$ws.Range("D24").Value = 0.03069205704343323
$ws.Range("E24").Value = 0.01926815323041731
$ws.Range("D25").Value = 0.01416151517098377
$ws.Range("E25").Value = 0.04448574969021069
$ws.Range("D26").Value = 0.01436390609412133
$ws.Range("E26").Value = 0.02858771786935965
$ws.Range("D27").Value = 0.03183308949256715
$ws.Range("E27").Value = -0.0199434038539279
$ws.Range("D28").Value = 0.02887131866168124
$ws.Range("E28").Value = 0.0676427683829488
$ws.Range("D29").Value = 0.02921565620528523
$ws.Range("E29").Value = 0.010931806350859
$ws.Range("D30").Value = 0.02863090632234732
$ws.Range("E30").Value = 0.03236197467975122
$ws.Range("D31").Value = 0.03398061707198633
$ws.Range("E31").Value = 0.01991668388074097
$ws.Range("D32").Value = 0.03121860491234747
$ws.Range("E32").Value = 0.0000343512761500353
$ws.Range("D33").Value = 0.02965670463316304
$ws.Range("E33").Value = -0.001452991452991492
$ws.Range("D34").Value = 0.03236460939283015
$ws.Range("E34").Value = 0.006458298190230582
$ws.Range("D35").Value = 0.03028142574851251
$ws.Range("E35").Value = 0.0003477051460361036
$ws.Range("D36").Value = 0.0313315967620182
$ws.Range("E36").Value = 0.003422739436181432
$ws.Range("D37").Value = 0.03322486830697555
$ws.Range("E37").Value = 0.008028169014084652
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0.00927927285893837

$ws.Protect("Lido2021")
